# Keep a single scenario for every triptype to run on jenkins.
$wb = $excel.ActiveWorkbook

# --- Sheet: Air_Mystifly_OneWay ---
$wsOneWay = $wb.Worksheets.Item("Air_Mystifly_OneWay")
$wsOneWay.Activate()

# Remove the child/infant scenario rows (3,4,5) - keep only the single "1 Adult" scenario.
$wsOneWay.Rows("3:5").Delete()

# Drop the trailing "|PAYNOW" step from the execution pipeline for the remaining scenario.
$wsOneWay.Range("B2").Value = "LOGIN|Search|AddToCart|CHECKOUTTRIP|ENTERPAXINFO|CONFIRMPAXINFO"
$wsOneWay.Range("B2").Select()

# --- Sheet: Air_Mystifly_RoundTrip ---
$wsRoundTrip = $wb.Worksheets.Item("Air_Mystifly_RoundTrip")
$wsRoundTrip.Activate()
$wsRoundTrip.Range("B2").Value = "LOGIN|Search|AddToCart|CHECKOUTTRIP|ENTERPAXINFO|CONFIRMPAXINFO"
$wsRoundTrip.Range("B7").Select()

# --- Sheet: Air_Mystifly_Multicity ---
$wsMulticity = $wb.Worksheets.Item("Air_Mystifly_Multicity")
$wsMulticity.Activate()
$wsMulticity.Range("B2").Value = "LOGIN|Search|AddToCart|CHECKOUTTRIP|ENTERPAXINFO|CONFIRMPAXINFO"
$wsMulticity.Range("B2").Select()

# Multicity is the scenario Jenkins should land on when the workbook opens.
$wsMulticity.Activate()
